$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1, columns B..I) with the new column labels
$ws.Range("B1").Value2 = "pess"
$ws.Range("C1").Value2 = "X0"
$ws.Range("D1").Value2 = "X20"
$ws.Range("E1").Value2 = "X40"
$ws.Range("F1").Value2 = "X60"
$ws.Range("G1").Value2 = "X80"
$ws.Range("H1").Value2 = "X100"
$ws.Range("I1").Value2 = "opt"

# Reorder the row labels in column A (rows 2..9) to match the new order
$ws.Range("A2").Value2 = "R...1.."
$ws.Range("A3").Value2 = "R...5.."
$ws.Range("A4").Value2 = "R...6.."
$ws.Range("A5").Value2 = "R...7.."
$ws.Range("A6").Value2 = "R...8.."
$ws.Range("A7").Value2 = "R...2.."
$ws.Range("A8").Value2 = "R...3.."
$ws.Range("A9").Value2 = "R...4.."

# Select A2 to match the saved cursor/selection state
$ws.Range("A2").Select()
